# Swap summary statistics between two pairs of rows in the single
# "BT_ID_loc_sum_pred_events_unfilt" table:
#   - Row for F59776 <-> Row for F59803  (Number of locations only)
#   - Row for F59783 <-> Row for F59809  (Number of locations, Number of
#     days tracked, Date sequence breaks in tracking, Number of
#     untracked days)

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Get-CellText($table, $row, $col) {
    $cell = $table.Cell($row, $col)
    $text = $cell.Range.Text
    # Word cell text ends with a cell-mark (bell + paragraph mark); strip
    # trailing control characters.
    return $text.TrimEnd([char]0x07, [char]0x0D, [char]0x0A)
}

function Set-CellText($table, $row, $col, $value) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $value
}

# Sanity-check the row identifiers (column 1 = ID) before touching anything,
# so the swap below is only ever applied to the intended rows.
$id1 = Get-CellText $t 46 1
$id2 = Get-CellText $t 55 1
$id3 = Get-CellText $t 52 1
$id4 = Get-CellText $t 60 1

if ($id1 -ne "F59776") { throw "Unexpected ID in row 46: $id1" }
if ($id2 -ne "F59803") { throw "Unexpected ID in row 55: $id2" }
if ($id3 -ne "F59783") { throw "Unexpected ID in row 52: $id3" }
if ($id4 -ne "F59809") { throw "Unexpected ID in row 60: $id4" }

# Row 46 = F59776 (Control, Roach); Row 55 = F59803 (Mix, Roach)
$r1 = 46
$r2 = 55
$col = 4   # "Number of locations"

$v1 = Get-CellText $t $r1 $col
$v2 = Get-CellText $t $r2 $col
Set-CellText $t $r1 $col $v2
Set-CellText $t $r2 $col $v1

# Row 52 = F59783 (Control, Roach); Row 60 = F59809 (Mix, Roach)
$r3 = 52
$r4 = 60
$cols = @(4, 5, 7, 8)   # locations, days tracked, date seq breaks, untracked days

foreach ($c in $cols) {
    $a = Get-CellText $t $r3 $c
    $b = Get-CellText $t $r4 $c
    Set-CellText $t $r3 $c $b
    Set-CellText $t $r4 $c $a
}

Write-Host "Swap complete"
